$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 602.5
$ws.Range("I42").Value = 112.888885
$ws.Range("K42").Value = 338.666655
$ws.Range("M42").Value = -108.666655

$ws.Range("H86").Value = 7126.1113
$ws.Range("I86").Value = 6524.5454
$ws.Range("J86").Value = 8071.4287
$ws.Range("K86").Value = 6524.5454
$ws.Range("L86").Value = 8071.4287
$ws.Range("M86").Value = -5401.5454
$ws.Range("N86").Value = -10317.4287

$ws.Range("H89").Value = 7126.1113
$ws.Range("I89").Value = 6524.5454
$ws.Range("J89").Value = 8071.4287
$ws.Range("K89").Value = 32622.727
$ws.Range("L89").Value = 40357.14350000001
$ws.Range("M89").Value = -27006.727
$ws.Range("N89").Value = -51589.14350000001

$ws.Range("H100").Value = 3748.5
$ws.Range("J100").Value = 4007.3333
$ws.Range("L100").Value = 4007.3333
$ws.Range("N100").Value = -5089.3333

$ws.Range("H113").Value = 100002296
$ws.Range("I113").Value = 50001300
$ws.Range("J113").Value = 133336296
$ws.Range("K113").Value = 50001300
$ws.Range("L113").Value = 133336296
$ws.Range("M113").Value = -49998046
$ws.Range("N113").Value = -133342804

$ws.Range("H135").Value = 6413.5
$ws.Range("I135").Value = 2104.8125
$ws.Range("K135").Value = 18943.3125
$ws.Range("M135").Value = -16408.3125

$ws.Range("H138").Value = 1451502.6
$ws.Range("I138").Value = 948.2
$ws.Range("J138").Value = 2567313.8
$ws.Range("K138").Value = 2844.6
$ws.Range("L138").Value = 7701941.399999999
$ws.Range("M138").Value = 2295.4
$ws.Range("N138").Value = -7712221.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11368465
$ws.Range("I32").Value = 11368465
$ws.Range("K32").Value = 11368465
$ws.Range("M32").Value = -11368178

$ws.Range("H57").Value = 12561.667
$ws.Range("I57").Value = 12561.667
$ws.Range("K57").Value = 12561.667
$ws.Range("M57").Value = -12077.667

$ws.Range("H74").Value = 6951364.5
$ws.Range("I74").Value = 9617034
$ws.Range("J74").Value = 20624.1
$ws.Range("K74").Value = 9617034
$ws.Range("L74").Value = 20624.1
$ws.Range("M74").Value = -9616160
$ws.Range("N74").Value = -22372.1

$ws.Range("H77").Value = 6951364.5
$ws.Range("I77").Value = 9617034
$ws.Range("J77").Value = 20624.1
$ws.Range("K77").Value = 48085170
$ws.Range("L77").Value = 103120.5
$ws.Range("M77").Value = -48080802
$ws.Range("N77").Value = -111856.5

$ws.Range("H88").Value = 1791.1538
$ws.Range("I88").Value = 1555.5714
$ws.Range("J88").Value = 2066
$ws.Range("K88").Value = 1555.5714
$ws.Range("L88").Value = 2066
$ws.Range("M88").Value = -1149.5714
$ws.Range("N88").Value = -2878

$ws.Range("H91").Value = 1791.1538
$ws.Range("I91").Value = 1555.5714
$ws.Range("J91").Value = 2066
$ws.Range("K91").Value = 1555.5714
$ws.Range("L91").Value = 2066
$ws.Range("M91").Value = -151.5714
$ws.Range("N91").Value = -4874

$ws.Range("H109").Value = 34995
$ws.Range("J109").Value = 34995
$ws.Range("L109").Value = 34995
$ws.Range("N109").Value = -37769

$ws.Range("H122").Value = 3979.8
$ws.Range("I122").Value = 2899.5
$ws.Range("J122").Value = 4249.875
$ws.Range("K122").Value = 8698.5
$ws.Range("L122").Value = 12749.625
$ws.Range("M122").Value = -6248.5
$ws.Range("N122").Value = -17649.625

$ws.Range("H132").Value = 7878.952
$ws.Range("I132").Value = 4549.5625
$ws.Range("K132").Value = 13648.6875
$ws.Range("M132").Value = -11118.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 66615
$ws.Range("J2").Value = 66615
$ws.Range("L2").Value = 66615
$ws.Range("N2").Value = -66841

$ws.Range("H86").Value = 1978.8235
$ws.Range("I86").Value = 2064.923
$ws.Range("J86").Value = 1699
$ws.Range("K86").Value = 2064.923
$ws.Range("L86").Value = 1699
$ws.Range("M86").Value = -941.9229999999998
$ws.Range("N86").Value = -3945

$ws.Range("H89").Value = 1978.8235
$ws.Range("I89").Value = 2064.923
$ws.Range("J89").Value = 1699
$ws.Range("K89").Value = 10324.615
$ws.Range("L89").Value = 8495
$ws.Range("M89").Value = -4708.614999999998
$ws.Range("N89").Value = -19727

$ws.Range("H94").Value = 2120.2144
$ws.Range("I94").Value = 1854.2222
$ws.Range("K94").Value = 1854.2222
$ws.Range("M94").Value = -1403.2222

$ws.Range("H107").Value = 1735.3334
$ws.Range("I107").Value = 1443.5714
$ws.Range("K107").Value = 1443.5714
$ws.Range("M107").Value = 476.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5877.7896
$ws.Range("I6").Value = 2042.5714
$ws.Range("J6").Value = 16616.4
$ws.Range("K6").Value = 2042.5714
$ws.Range("L6").Value = 16616.4
$ws.Range("M6").Value = -1929.5714
$ws.Range("N6").Value = -16842.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 10008
$ws.Range("I17").Value = 10008
$ws.Range("K17").Value = 10008
$ws.Range("M17").Value = -9834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 7430.278
$ws.Range("J25").Value = 7749.1177
$ws.Range("L25").Value = 7749.1177
$ws.Range("N25").Value = -8097.1177

$ws.Range("H51").Value = 40285.715
$ws.Range("I51").Value = 18250
$ws.Range("J51").Value = 69666.664
$ws.Range("K51").Value = 18250
$ws.Range("L51").Value = 69666.664
$ws.Range("M51").Value = -17514
$ws.Range("N51").Value = -71138.664

$ws.Range("H61").Value = 40285.715
$ws.Range("I61").Value = 18250
$ws.Range("J61").Value = 69666.664
$ws.Range("K61").Value = 18250
$ws.Range("L61").Value = 69666.664
$ws.Range("M61").Value = -17902
$ws.Range("N61").Value = -70362.664

$ws.Range("H64").Value = 109000
$ws.Range("J64").Value = 109000
$ws.Range("L64").Value = 109000
$ws.Range("N64").Value = -109496

$ws.Range("H67").Value = 109000
$ws.Range("J67").Value = 109000
$ws.Range("L67").Value = 109000
$ws.Range("N67").Value = -110716

$ws.Range("H132").Value = 2043.9412
$ws.Range("I132").Value = 1859.1875
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5577.5625
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3047.5625
$ws.Range("N132").Value = -20060

$ws.Range("H134").Value = 3340668.2
$ws.Range("I134").Value = 10000000
$ws.Range("J134").Value = 11002.5
$ws.Range("K134").Value = 30000000
$ws.Range("L134").Value = 33007.5
$ws.Range("M134").Value = -29997465
$ws.Range("N134").Value = -38077.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 3334666.8
$ws.Range("I92").Value = 3334666.8
$ws.Range("K92").Value = 10004000.4
$ws.Range("M92").Value = -10002752.4

$ws.Range("H113").Value = 1137.7273
$ws.Range("J113").Value = 1509
$ws.Range("L113").Value = 4527
$ws.Range("N113").Value = -8867

$ws.Range("H117").Value = 233
$ws.Range("I117").Value = 233
$ws.Range("K117").Value = 699
$ws.Range("M117").Value = 2743

$ws.Range("H132").Value = 2586.0667
$ws.Range("I132").Value = 2390.0908
$ws.Range("K132").Value = 21510.8172
$ws.Range("M132").Value = -18980.8172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1443.1428
$ws.Range("I97").Value = 1501.6
$ws.Range("K97").Value = 1501.6
$ws.Range("M97").Value = -1005.6

$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344

$ws.Range("H113").Value = 3506.5334
$ws.Range("I113").Value = 2469.8
$ws.Range("J113").Value = 4024.9
$ws.Range("K113").Value = 2469.8
$ws.Range("L113").Value = 4024.9
$ws.Range("M113").Value = -299.8000000000002
$ws.Range("N113").Value = -8364.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1481.875
$ws.Range("I16").Value = 1586.7142
$ws.Range("K16").Value = 1586.7142
$ws.Range("M16").Value = -1416.7142

$ws.Range("H82").Value = 715.5
$ws.Range("I82").Value = 703.5714
$ws.Range("K82").Value = 703.5714
$ws.Range("M82").Value = -342.5714

$ws.Range("H85").Value = 715.5
$ws.Range("I85").Value = 703.5714
$ws.Range("K85").Value = 703.5714
$ws.Range("M85").Value = 544.4286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4553132
$ws.Range("J62").Value = 33340000
$ws.Range("L62").Value = 33340000
$ws.Range("N62").Value = -33341248

$ws.Range("H65").Value = 4553132
$ws.Range("J65").Value = 33340000
$ws.Range("L65").Value = 166700000
$ws.Range("N65").Value = -166706240

$ws.Range("H113").Value = 1205.1538
$ws.Range("I113").Value = 1242.4445
$ws.Range("K113").Value = 3727.3335
$ws.Range("M113").Value = -1557.3335

$ws.Range("H122").Value = 4366.8184
$ws.Range("I122").Value = 2717.9565
$ws.Range("K122").Value = 8153.869499999999
$ws.Range("M122").Value = -5703.869499999999
